$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates -----------------------------------------------------
# A3: "Testing" -> "Multi line items"
$ws.Range("A3").Value = "Multi line items"
# B3 was empty -> now holds the placeholder quoted text
$ws.Range("B3").Value = [char]8220 + " " + [char]8221

# H3 becomes the new "clearCartItems" action (new Arial/black font),
# I3 and J3 (old productDetailPage / cartCheck) are cleared out.
$ws.Range("H3").Value = "clearCartItems"
$ws.Range("H3").Font.Name = "Arial"
$ws.Range("H3").Font.Size = 10
$ws.Range("H3").Font.Color = 0
$ws.Range("I3").ClearContents()
$ws.Range("J3").Clear()

# --- Row 5 update --------------------------------------------------------
$ws.Range("C5").Value = "YES"

# --- New row 10: clearCartItems test case --------------------------------
$ws.Range("B10").Value = "checking cartitemclear"
$ws.Range("C10").Value = "NO"
$ws.Range("D10").Value = "clearCartItems"

# Re-use the exact same style as H3 (same font) by copying H3's format
# onto D10, so both cells share one style entry, same as the source file.
$ws.Range("H3").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Cosmetic layout tweaks to mirror the authored workbook --------------
$ws.Rows.Item(6).RowHeight = 35.05

$ws.Columns.Item(1).ColumnWidth = 13.748299319727867
$ws.Columns.Item(2).ColumnWidth = 23.059523809523768
$ws.Columns.Item(3).ColumnWidth = 4.702380952380957
$ws.Columns.Item(4).ColumnWidth = 23.059523809523768
$ws.Columns.Item(5).ColumnWidth = 21.17176870748297
$ws.Columns.Item(6).ColumnWidth = 20.493197278911566
$ws.Columns.Item(7).ColumnWidth = 16.17687074829937
$ws.Columns.Item(8).ColumnWidth = 21.845238095238066
$ws.Columns.Item(9).ColumnWidth = 18.73809523809527

$ws.Range("C5").Select()

Write-Output "done"
